$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Re-apply the AutoFilter: switch the filtered column from "association" (col B,
# field 2) value "38" to "nais1" (col J, field 10) value "18M". Clear the old
# filter first so only the new filterColumn remains, and use the "values"
# operator (xlFilterValues = 7) so the XML is written as a discrete <filters>
# list (matching Excel's own output) rather than a <customFilters> entry.
$ws.AutoFilterMode = $False
$rng = $ws.Range("A1:N928")
$rng.AutoFilter(10, @("18M"), 7)

# Row 465's "nais1" text has a trailing space ("18M "), so Excel's own
# whitespace-insensitive match still shows it for the "18M" filter even
# though a strict equality wouldn't; make sure it stays visible.
$ws.Rows.Item(465).Hidden = $False

# Update the active selection / active cell shown in the sheet view.
$ws.Activate()
$ws.Range("J437").Select()
